# ------------------------------------------------------------------
# Layout edit for "learn more about subsidies" page
#   1. Re-brand all "Avenir Medium" / "Avenir Heavy" runs to "Avenir Roman"
#   2. Remove the two trailing blank paragraphs and the tab/bookmark-only
#      paragraph at the end of the story
#   3. Insert "each year" into the closing sentence and relocate the
#      "_GoBack" bookmark so it now sits inside that sentence
# ------------------------------------------------------------------

$d = $word.ActiveDocument

# --- 1. Normalize all fonts to "Avenir Roman" --------------------------
$paraCount = $d.Paragraphs.Count
for ($i = 1; $i -le $paraCount; $i++) {
    $d.Paragraphs($i).Range.Font.Name = "Avenir Roman"
}

# --- 2. Drop the two empty paragraphs + the tab/bookmark paragraph -----
# (paragraph 4 is "This could mean ... !", paragraphs 5 & 6 are blank,
#  paragraph 7 holds just a tab + the _GoBack bookmark)
$firstBlank = $d.Paragraphs(5)
$tabPara = $d.Paragraphs(7)
$trailing = $d.Range($firstBlank.Range.Start, $tabPara.Range.End)
$trailing.Delete()

# --- 3. Split the closing sentence, add "each year" and re-home the
#        _GoBack bookmark inside it ------------------------------------
$closing = $d.Paragraphs(4)
$paraStart = $closing.Range.Start
$paraText = $closing.Range.Text

# Locate "the water" so the bookmark can be re-inserted between
# "th" and "e water", exactly where Word itself split the run.
$wordOffset = $paraText.IndexOf("the water")
$bookmarkPos = $paraStart + $wordOffset + 2

$bookmarkRange = $d.Range($bookmarkPos, $bookmarkPos)
$d.Bookmarks.Add("_GoBack", $bookmarkRange)

# Insert " each year" right before the closing "!"
$closingNow = $d.Paragraphs(4)
$paraEnd = $closingNow.Range.End
$beforeBang = $paraEnd - 2   # skip the paragraph mark, then the "!" char
$insertion = " each year"
$insertRange = $d.Range($beforeBang, $beforeBang)
$insertRange.InsertBefore($insertion)

# --- 4. Break the tail of the sentence into its own runs (matching how
#        Word itself would leave separate runs after an in-place edit),
#        without altering any visible formatting. Toggling Bold on/off
#        forces a run boundary while leaving rPr content unchanged
#        (the run already carries <w:bCs/> rather than <w:b/>).
$finalPara = $d.Paragraphs(4)
$finalEnd = $finalPara.Range.End

$bangRange = $d.Range($finalEnd - 2, $finalEnd - 1)
$bangRange.Bold = 1
$bangRange.Bold = 0

$eachYearRange = $d.Range($finalEnd - 2 - $insertion.Length, $finalEnd - 2)
$eachYearRange.Bold = 1
$eachYearRange.Bold = 0
